# --- Update the confidential disclosure date (shared string behind A80) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow writes, then restore protection at the end.
$ws.Unprotect("")

$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# --- Update Weight (D) / Percent Change (E) columns for rows 2-77 ---
$data = New-Object 'object[,]' 76,2
$data[0,0] = 0.07435950159377423
$data[0,1] = -0.01227448719004864
$data[1,0] = 0.04521817279914706
$data[1,1] = -0.006645576822413068
$data[2,0] = 0.03573686160655758
$data[2,1] = -0.01441081448733228
$data[3,0] = 0.03363884431865401
$data[3,1] = -0.01216726492937381
$data[4,0] = 0.03142497836989359
$data[4,1] = 0.01178936337437775
$data[5,0] = 0.03007639998659409
$data[5,1] = 0.0003323899324954027
$data[6,0] = 0.03092002718483897
$data[6,1] = -0.006203324500120511
$data[7,0] = 0.02917220216871339
$data[7,1] = -0.002434124101692303
$data[8,0] = 0.02628902107190601
$data[8,1] = -0.003422487007225361
$data[9,0] = 0.02831500048694172
$data[9,1] = -0.01541033655593516
$data[10,0] = 0.02355246399933031
$data[10,1] = 0.003678658371652821
$data[11,0] = 0.02365517922486161
$data[11,1] = 0.01774993474288711
$data[12,0] = 0.01973367265164339
$data[12,1] = -0.00518520725929017
$data[13,0] = 0.01864320547677583
$data[13,1] = -0.004289859948689911
$data[14,0] = 0.02136221863201363
$data[14,1] = -0.01428027418126432
$data[15,0] = 0.0187531734949496
$data[15,1] = -0.006804712080193065
$data[16,0] = 0.01845081045032742
$data[16,1] = -0.01567571300324566
$data[17,0] = 0.01547099079209459
$data[17,1] = -0.01236933797909412
$data[18,0] = 0.01412192235523817
$data[18,1] = 0.0005621662138106664
$data[19,0] = 0.0156769112967141
$data[19,1] = -0.009696719620383765
$data[20,0] = 0.0142276269074667
$data[20,1] = -0.01297842440274455
$data[21,0] = 0.0129888695262739
$data[21,1] = -0.003870967741935405
$data[22,0] = 0.01527555743359325
$data[22,1] = 0.005562827225130906
$data[23,0] = 0.01486077610301267
$data[23,1] = -0.01450957632037142
$data[24,0] = 0.01188144649833677
$data[24,1] = -0.0005774338838201443
$data[25,0] = 0.0121601889615152
$data[25,1] = 0.008946562424438032
$data[26,0] = 0.01232729722442633
$data[26,1] = -0.01131782945736426
$data[27,0] = 0.01190026455492265
$data[27,1] = -0.003805037144410162
$data[28,0] = 0.01117890571913034
$data[28,1] = 0.003980431008784979
$data[29,0] = 0.01276099462230418
$data[29,1] = 0.004608294930875667
$data[30,0] = 0.01332847664122232
$data[30,1] = 0.008669755129053769
$data[31,0] = 0.01123619297993477
$data[31,1] = -0.01548291427699133
$data[32,0] = 0.01182023880907694
$data[32,1] = -0.004079551249362434
$data[33,0] = 0.008986945164364428
$data[33,1] = 0.03980107641217767
$data[34,0] = 0.01134743513735655
$data[34,1] = -0.01299907149489343
$data[35,0] = 0.01111343456392528
$data[35,1] = -0.008907311050357247
$data[36,0] = 0.0101845380467756
$data[36,1] = 0.00801154817755334
$data[37,0] = 0.009264364482939175
$data[37,1] = -0.03482713384960434
$data[38,0] = 0.009583879402053705
$data[38,1] = -0.02028961793340422
$data[39,0] = 0.009006106258440163
$data[39,1] = -0.008836748685914553
$data[40,0] = 0.009158414903932181
$data[40,1] = -0.004045246834969096
$data[41,0] = 0.009976510311861312
$data[41,1] = -0.01711366538952752
$data[42,0] = 0.009760102661123617
$data[42,1] = -0.008154084071418488
$data[43,0] = 0.009431374735138916
$data[43,1] = -0.009259259259259078
$data[44,0] = 0.009690711077463162
$data[44,1] = -0.0139268159475695
$data[45,0] = 0.008745691798290682
$data[45,1] = 0.0252824098977944
$data[46,0] = 0.007305326383779298
$data[46,1] = -0.004507888805409532
$data[47,0] = 0.008270829901631206
$data[47,1] = 0.00763150722267647
$data[48,0] = 0.008021931700069443
$data[48,1] = 0.01044625675799504
$data[49,0] = 0.007897409091255024
$data[49,1] = 0.01242289983494071
$data[50,0] = 0.007852128142595236
$data[50,1] = -0.002883355176933056
$data[51,0] = 0.00695680029409486
$data[51,1] = 0.03338968723584124
$data[52,0] = 0.007427839773010329
$data[52,1] = -0.01134774232707891
$data[53,0] = 0.006713385791835079
$data[53,1] = -0.0008606285435229788
$data[54,0] = 0.006591808404897774
$data[54,1] = 0.002973712382538407
$data[55,0] = 0.006835659054823217
$data[55,1] = 0.0003441156228494169
$data[56,0] = 0.006314046048833181
$data[56,1] = 0.01927912824811395
$data[57,0] = 0.005618611046202156
$data[57,1] = 0.007675333397294315
$data[58,0] = 0.006459150907038414
$data[58,1] = 0.01559121429384303
$data[59,0] = 0.005399410090189996
$data[59,1] = 0.01207115628970779
$data[60,0] = 0.005915926539188023
$data[60,1] = -0.01418157720344604
$data[61,0] = 0.005445867167386403
$data[61,1] = -0.004319343459794123
$data[62,0] = 0.005023637022740577
$data[62,1] = -0.01186202590916174
$data[63,0] = 0.004844473442329118
$data[63,1] = -0.007890264627336663
$data[64,0] = 0.004546030826162327
$data[64,1] = -0.02858806028070637
$data[65,0] = 0.004492811009880367
$data[65,1] = -0.005890052356020914
$data[66,0] = 0.003464727652809248
$data[66,1] = 0.00369160266474311
$data[67,0] = 0.004104982624930616
$data[67,1] = 0.02399541580116016
$data[68,0] = 0.003647272602763484
$data[68,1] = 0.02276086313922554
$data[69,0] = 0.003098657645786927
$data[69,1] = -0.003653271338425856
$data[70,0] = 0.002538134387377114
$data[70,1] = -0.005908134303863388
$data[71,0] = 0.002518630255811534
$data[71,1] = -0.001089600155657289
$data[72,0] = 0.002299919353356297
$data[72,1] = -0.008224665473451021
$data[73,0] = 0.00186828018041753
$data[73,1] = -0.003777148253068963
$data[74,0] = 0.001758410172955141
$data[74,1] = 0.03979711275848596
$data[75,0] = 1
$data[75,1] = -0.003891559400354683
$ws.Range("D2:E77").Value = $data

# --- Restore sheet protection ---
$ws.Protect("")
